# Apply updated crypto price/volume data to the worksheet.
# Numeric-looking Price strings are entered with a leading apostrophe so Excel
# keeps them as text (matching the source data) instead of auto-converting them
# to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.293.35"
$ws.Range("E2").Value = "  +4.91%  "

$ws.Range("D3").Value = "1.701.91"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("D4").Value = "'0.9953"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "'240.41"
$ws.Range("E5").Value = "  +3.03%  "

$ws.Range("D6").Value = "'0.9972"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").Value = "'0.4693"
$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("D8").Value = "'0.2641"
$ws.Range("E8").Value = "  +3.02%  "

$ws.Range("D9").Value = "'0.06202"
$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").Value = "1.689.38"
$ws.Range("E10").Value = "  +2.64%  "

$ws.Range("D11").Value = "'0.07068"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("E12").Value = "  +6.23%  "

$ws.Range("D13").Value = "'4.419"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").Value = "'0.5888"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").Value = "'76.12"
$ws.Range("E15").Value = "  +3.69%  "

$ws.Range("D16").Value = "'0.9978"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").Value = "'0.9978"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "26.248.00"
$ws.Range("E18").Value = "  +4.73%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000006806"
$ws.Range("E19").Value = "  +3.39%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'11.60"
$ws.Range("E20").Value = "  +2.97%  "

$ws.Range("D21").Value = "1.900.71"
$ws.Range("E21").Value = "  +2.75%  "

$ws.Range("D22").Value = "'4.550"
$ws.Range("E22").Value = "  +5.88%  "

$ws.Range("D23").Value = "'8.804"
$ws.Range("E23").Value = "  +4.05%  "

$ws.Range("D24").Value = "'5.320"
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("D25").Value = "'134.93"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").Value = "'15.15"
$ws.Range("E26").Value = "  +1.70%  "

$ws.Range("D27").Value = "'1.397"
$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").Value = "'1.748"
$ws.Range("E28").Value = "  +6.87%  "

$ws.Range("D29").Value = "'106.09"
$ws.Range("E29").Value = "  +2.57%  "

$ws.Range("D30").Value = "'4.002"
$ws.Range("E30").Value = "  +3.03%  "

$ws.Range("D31").Value = "'3.695"
$ws.Range("E31").Value = "  +4.96%  "

$ws.Range("D32").Value = "'0.07775"
$ws.Range("E32").Value = "  +2.63%  "

$ws.Range("D33").Value = "'0.04394"
$ws.Range("E33").Value = "  +3.71%  "

$ws.Range("D34").Value = "'2.597"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6192"
$ws.Range("E35").Value = "  +4.94%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9712"
$ws.Range("E36").Value = "  +4.18%  "

$ws.Range("D37").Value = "'0.9198"
$ws.Range("E37").Value = "  +6.85%  "

$ws.Range("D38").Value = "'110.76"
$ws.Range("E38").Value = "  +12.88%  "

$ws.Range("E39").Value = "  -7.92%  "

$ws.Range("D40").Value = "'0.9992"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'1.901"
$ws.Range("E41").Value = "  +6.93%  "

$ws.Range("D42").Value = "'0.01467"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").Value = "'0.3787"
$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("D44").Value = "'5.116"
$ws.Range("E44").Value = "  +10.76%  "

$ws.Range("D45").Value = "'0.1135"
$ws.Range("E45").Value = "  +3.71%  "

$ws.Range("D46").Value = "'6.243"
$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("D47").Value = "'0.05318"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").Value = "'30.88"
$ws.Range("E48").Value = "  +7.23%  "

$ws.Range("D49").Value = "'7.660"
$ws.Range("E49").Value = "  +7.10%  "

$ws.Range("D50").Value = "'1.219"
$ws.Range("E50").Value = "  +2.29%  "

$ws.Range("D51").Value = "'0.3365"
$ws.Range("E51").Value = "  +2.85%  "
